$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G2").Value = 128.4548946666667
$ws.Range("H2").Value = 385.364684
$ws.Range("I2").Value = 0.2815548034715028
$ws.Range("J2").Value = 0.2815548034715028
$ws.Range("M2").Value = 5.455768666666667
$ws.Range("N2").Value = 16.367306
$ws.Range("O2").Value = 0.0824390136851795
$ws.Range("P2").Value = 0.0824390136851795
$ws.Range("Q2").Value = 700.8201894023671
$ws.Range("R2").Value = 6307.381704621303
$ws.Range("S2").Value = 0.02321110029651524
$ws.Range("T2").Value = 0.02321110029651524
$ws.Range("G3").Value = 128.4548946666667
$ws.Range("H3").Value = 385.364684
$ws.Range("I3").Value = 0.2815548034715028
$ws.Range("J3").Value = 0.2815548034715028
$ws.Range("O3").Value = 0.03476630532942922
$ws.Range("P3").Value = 0.03476630532942922
$ws.Range("Q3").Value = 295.5509484724867
$ws.Range("R3").Value = 2659.95853625238
$ws.Range("S3").Value = 0.009788620264457705
$ws.Range("T3").Value = 0.009788620264457705
$ws.Range("G4").Value = 128.4548946666667
$ws.Range("H4").Value = 385.364684
$ws.Range("I4").Value = 0.2815548034715028
$ws.Range("J4").Value = 0.2815548034715028
$ws.Range("M4").Value = 29.166511
$ws.Range("N4").Value = 87.499533
$ws.Range("O4").Value = 0.4407185396566677
$ws.Range("P4").Value = 0.4407185396566677
$ws.Range("Q4").Value = 3746.581098299174
$ws.Range("R4").Value = 33719.22988469257
$ws.Range("S4").Value = 0.1240864218192808
$ws.Range("T4").Value = 0.1240864218192808
$ws.Range("G5").Value = 128.4548946666667
$ws.Range("H5").Value = 385.364684
$ws.Range("I5").Value = 0.2815548034715028
$ws.Range("J5").Value = 0.2815548034715028
$ws.Range("M5").Value = 3.497096
$ws.Range("N5").Value = 10.491288
$ws.Range("O5").Value = 0.05284262633124592
$ws.Range("P5").Value = 0.05284262633124593
$ws.Range("Q5").Value = 449.2190983192212
$ws.Range("R5").Value = 4042.971884872992
$ws.Range("S5").Value = 0.014878095271612
$ws.Range("T5").Value = 0.01487809527161201
$ws.Range("G6").Value = 128.4548946666667
$ws.Range("H6").Value = 385.364684
$ws.Range("I6").Value = 0.2815548034715028
$ws.Range("J6").Value = 0.2815548034715028
$ws.Range("M6").Value = 25.75926033333333
$ws.Range("N6").Value = 77.277781
$ws.Range("O6").Value = 0.3892335149974776
$ws.Range("P6").Value = 0.3892335149974776
$ws.Range("Q6").Value = 3308.903072809578
$ws.Range("R6").Value = 29780.12765528621
$ws.Range("S6").Value = 0.109590565819637
$ws.Range("T6").Value = 0.109590565819637
$ws.Range("H7").Value = 457.183265
$ws.Range("I7").Value = 0.3340268313936494
$ws.Range("J7").Value = 0.3340268313936494
$ws.Range("M7").Value = 5.455768666666667
$ws.Range("N7").Value = 16.367306
$ws.Range("O7").Value = 0.0824390136851795
$ws.Range("P7").Value = 0.0824390136851795
$ws.Range("Q7").Value = 831.4287107037877
$ws.Range("R7").Value = 7482.85839633409
$ws.Range("S7").Value = 0.02753684252447821
$ws.Range("T7").Value = 0.02753684252447821
$ws.Range("H8").Value = 457.183265
$ws.Range("I8").Value = 0.3340268313936494
$ws.Range("J8").Value = 0.3340268313936494
$ws.Range("O8").Value = 0.03476630532942922
$ws.Range("P8").Value = 0.03476630532942922
$ws.Range("S8").Value = 0.01161287880845339
$ws.Range("T8").Value = 0.01161287880845339
$ws.Range("H9").Value = 457.183265
$ws.Range("I9").Value = 0.3340268313936494
$ws.Range("J9").Value = 0.3340268313936494
$ws.Range("M9").Value = 29.166511
$ws.Range("N9").Value = 87.499533
$ws.Range("O9").Value = 0.4407185396566677
$ws.Range("P9").Value = 0.4407185396566677
$ws.Range("Q9").Value = 4444.813575879472
$ws.Range("R9").Value = 40003.32218291525
$ws.Range("S9").Value = 0.1472118173379532
$ws.Range("T9").Value = 0.1472118173379532
$ws.Range("H10").Value = 457.183265
$ws.Range("I10").Value = 0.3340268313936494
$ws.Range("J10").Value = 0.3340268313936494
$ws.Range("M10").Value = 3.497096
$ws.Range("N10").Value = 10.491288
$ws.Range("O10").Value = 0.05284262633124592
$ws.Range("P10").Value = 0.05284262633124593
$ws.Range("Q10").Value = 532.9379224328133
$ws.Range("R10").Value = 4796.44130189532
$ws.Range("S10").Value = 0.0176508550359447
$ws.Range("T10").Value = 0.0176508550359447
$ws.Range("H11").Value = 457.183265
$ws.Range("I11").Value = 0.3340268313936494
$ws.Range("J11").Value = 0.3340268313936494
$ws.Range("M11").Value = 25.75926033333333
$ws.Range("N11").Value = 77.277781
$ws.Range("O11").Value = 0.3892335149974776
$ws.Range("P11").Value = 0.3892335149974776
$ws.Range("Q11").Value = 3925.56758105944
$ws.Range("R11").Value = 35330.10822953497
$ws.Range("S11").Value = 0.13001443768682
$ws.Range("T11").Value = 0.13001443768682
$ws.Range("G12").Value = 70.798157
$ws.Range("H12").Value = 212.394471
$ws.Range("I12").Value = 0.1551794599342134
$ws.Range("J12").Value = 0.1551794599342134
$ws.Range("M12").Value = 5.455768666666667
$ws.Range("N12").Value = 16.367306
$ws.Range("O12").Value = 0.0824390136851795
$ws.Range("P12").Value = 0.0824390136851795
$ws.Range("Q12").Value = 386.2583666183473
$ws.Range("R12").Value = 3476.325299565126
$ws.Range("S12").Value = 0.01279284162117538
$ws.Range("T12").Value = 0.01279284162117538
$ws.Range("G13").Value = 70.798157
$ws.Range("H13").Value = 212.394471
$ws.Range("I13").Value = 0.1551794599342134
$ws.Range("J13").Value = 0.1551794599342134
$ws.Range("O13").Value = 0.03476630532942922
$ws.Range("P13").Value = 0.03476630532942922
$ws.Range("Q13").Value = 162.893461597955
$ws.Range("R13").Value = 1466.041154381595
$ws.Range("S13").Value = 0.005395016484928791
$ws.Range("T13").Value = 0.005395016484928791
$ws.Range("G14").Value = 70.798157
$ws.Range("H14").Value = 212.394471
$ws.Range("I14").Value = 0.1551794599342134
$ws.Range("J14").Value = 0.1551794599342134
$ws.Range("M14").Value = 29.166511
$ws.Range("N14").Value = 87.499533
$ws.Range("O14").Value = 0.4407185396566677
$ws.Range("P14").Value = 0.4407185396566677
$ws.Range("Q14").Value = 2064.935224920227
$ws.Range("R14").Value = 18584.41702428204
$ws.Range("S14").Value = 0.06839046496691691
$ws.Range("T14").Value = 0.06839046496691691
$ws.Range("G15").Value = 70.798157
$ws.Range("H15").Value = 212.394471
$ws.Range("I15").Value = 0.1551794599342134
$ws.Range("J15").Value = 0.1551794599342134
$ws.Range("M15").Value = 3.497096
$ws.Range("N15").Value = 10.491288
$ws.Range("O15").Value = 0.05284262633124592
$ws.Range("P15").Value = 0.05284262633124593
$ws.Range("Q15").Value = 247.587951652072
$ws.Range("R15").Value = 2228.291564868648
$ws.Range("S15").Value = 0.008200090215588186
$ws.Range("T15").Value = 0.008200090215588188
$ws.Range("G16").Value = 70.798157
$ws.Range("H16").Value = 212.394471
$ws.Range("I16").Value = 0.1551794599342134
$ws.Range("J16").Value = 0.1551794599342134
$ws.Range("M16").Value = 25.75926033333333
$ws.Range("N16").Value = 77.277781
$ws.Range("O16").Value = 0.3892335149974776
$ws.Range("P16").Value = 0.3892335149974776
$ws.Range("Q16").Value = 1823.708157283206
$ws.Range("R16").Value = 16413.37341554885
$ws.Range("S16").Value = 0.06040104664560413
$ws.Range("T16").Value = 0.06040104664560413
$ws.Range("G17").Value = 20.703408
$ws.Range("H17").Value = 62.110224
$ws.Range("I17").Value = 0.04537891674549766
$ws.Range("J17").Value = 0.04537891674549767
$ws.Range("M17").Value = 5.455768666666667
$ws.Range("N17").Value = 16.367306
$ws.Range("O17").Value = 0.0824390136851795
$ws.Range("P17").Value = 0.0824390136851795
$ws.Range("Q17").Value = 112.953004659616
$ws.Range("R17").Value = 1016.577041936544
$ws.Range("S17").Value = 0.003740993138600703
$ws.Range("T17").Value = 0.003740993138600704
$ws.Range("G18").Value = 20.703408
$ws.Range("H18").Value = 62.110224
$ws.Range("I18").Value = 0.04537891674549766
$ws.Range("J18").Value = 0.04537891674549767
$ws.Range("O18").Value = 0.03476630532942922
$ws.Range("P18").Value = 0.03476630532942922
$ws.Range("Q18").Value = 47.63471167752
$ws.Range("R18").Value = 428.71240509768
$ws.Range("S18").Value = 0.00157765727509272
$ws.Range("T18").Value = 0.001577657275092721
$ws.Range("G19").Value = 20.703408
$ws.Range("H19").Value = 62.110224
$ws.Range("I19").Value = 0.04537891674549766
$ws.Range("J19").Value = 0.04537891674549767
$ws.Range("M19").Value = 29.166511
$ws.Range("N19").Value = 87.499533
$ws.Range("O19").Value = 0.4407185396566677
$ws.Range("P19").Value = 0.4407185396566677
$ws.Range("Q19").Value = 603.846177169488
$ws.Range("R19").Value = 5434.615594525392
$ws.Range("S19").Value = 0.01999932991927723
$ws.Range("T19").Value = 0.01999932991927724
$ws.Range("G20").Value = 20.703408
$ws.Range("H20").Value = 62.110224
$ws.Range("I20").Value = 0.04537891674549766
$ws.Range("J20").Value = 0.04537891674549767
$ws.Range("M20").Value = 3.497096
$ws.Range("N20").Value = 10.491288
$ws.Range("O20").Value = 0.05284262633124592
$ws.Range("P20").Value = 0.05284262633124593
$ws.Range("Q20").Value = 72.40180530316799
$ws.Range("R20").Value = 651.616247728512
$ws.Range("S20").Value = 0.002397941140899051
$ws.Range("T20").Value = 0.002397941140899052
$ws.Range("G21").Value = 20.703408
$ws.Range("H21").Value = 62.110224
$ws.Range("I21").Value = 0.04537891674549766
$ws.Range("J21").Value = 0.04537891674549767
$ws.Range("M21").Value = 25.75926033333333
$ws.Range("N21").Value = 77.277781
$ws.Range("O21").Value = 0.3892335149974776
$ws.Range("P21").Value = 0.3892335149974776
$ws.Range("Q21").Value = 533.3044764592159
$ws.Range("R21").Value = 4799.740288132944
$ws.Range("S21").Value = 0.01766299527162795
$ws.Range("T21").Value = 0.01766299527162795
$ws.Range("G22").Value = 83.88319133333333
$ws.Range("H22").Value = 251.649574
$ws.Range("I22").Value = 0.1838599884551367
$ws.Range("J22").Value = 0.1838599884551367
$ws.Range("M22").Value = 5.455768666666667
$ws.Range("N22").Value = 16.367306
$ws.Range("O22").Value = 0.0824390136851795
$ws.Range("P22").Value = 0.0824390136851795
$ws.Range("Q22").Value = 457.6472869364049
$ws.Range("R22").Value = 4118.825582427644
$ws.Range("S22").Value = 0.01515723610440996
$ws.Range("T22").Value = 0.01515723610440996
$ws.Range("G23").Value = 83.88319133333333
$ws.Range("H23").Value = 251.649574
$ws.Range("I23").Value = 0.1838599884551367
$ws.Range("J23").Value = 0.1838599884551367
$ws.Range("O23").Value = 0.03476630532942922
$ws.Range("P23").Value = 0.03476630532942922
$ws.Range("Q23").Value = 192.9997048676033
$ws.Range("R23").Value = 1736.99734380843
$ws.Range("S23").Value = 0.006392132496496614
$ws.Range("T23").Value = 0.006392132496496614
$ws.Range("G24").Value = 83.88319133333333
$ws.Range("H24").Value = 251.649574
$ws.Range("I24").Value = 0.1838599884551367
$ws.Range("J24").Value = 0.1838599884551367
$ws.Range("M24").Value = 29.166511
$ws.Range("N24").Value = 87.499533
$ws.Range("O24").Value = 0.4407185396566677
$ws.Range("P24").Value = 0.4407185396566677
$ws.Range("Q24").Value = 2446.580022738771
$ws.Range("R24").Value = 22019.22020464894
$ws.Range("S24").Value = 0.08103050561323963
$ws.Range("T24").Value = 0.08103050561323963
$ws.Range("G25").Value = 83.88319133333333
$ws.Range("H25").Value = 251.649574
$ws.Range("I25").Value = 0.1838599884551367
$ws.Range("J25").Value = 0.1838599884551367
$ws.Range("M25").Value = 3.497096
$ws.Range("N25").Value = 10.491288
$ws.Range("O25").Value = 0.05284262633124592
$ws.Range("P25").Value = 0.05284262633124593
$ws.Range("Q25").Value = 293.3475728790346
$ws.Range("R25").Value = 2640.128155911312
$ws.Range("S25").Value = 0.009715644667201977
$ws.Range("T25").Value = 0.009715644667201978
$ws.Range("G26").Value = 83.88319133333333
$ws.Range("H26").Value = 251.649574
$ws.Range("I26").Value = 0.1838599884551367
$ws.Range("J26").Value = 0.1838599884551367
$ws.Range("M26").Value = 25.75926033333333
$ws.Range("N26").Value = 77.277781
$ws.Range("O26").Value = 0.3892335149974776
$ws.Range("P26").Value = 0.3892335149974776
$ws.Range("Q26").Value = 2160.768963146144
$ws.Range("R26").Value = 19446.9206683153
$ws.Range("S26").Value = 0.07156446957378851
$ws.Range("T26").Value = 0.07156446957378851
